# Update countries & provincias Spain
#
# Mirrors the source commit: bump the "Datos actualizados..." timestamp,
# refresh the case counters for a batch of already-listed countries, and
# splice in fresh rows for Albania (rank 98) and Siria (rank 150) -
# pushing the countries that used to sit at those ranking spots down by
# one row each (Paraguay; and Andorra/Chad/Malta/Jamaica respectively).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 5 de Agosto de 2020 a las 18:21"

# --- Per-row data: country name (column A) + B..H counters -----------
$rows = @{
    4   = @("Estados Unidos",          4927902, 9482,  2485771, 2281601, 0, 240, 160530)
    5   = @("Brasil",                  2817473, 9397,  1970767, 750380,  0, 230, 96326)
    6   = @("India",                   1959468, 52855, 1324568, 594171,  0, 909, 40729)
    11  = @("Chile",                   364723,  1761,  338291,  16640,   0, 47,  9792)
    12  = @("España",                  352847,  2953,  0,       0,       0, 1,   28499)
    15  = @("Reino Unido",              307184,  891,   0,       0,       0, 65,  46364)
    18  = @("Italia",                  248803,  384,   200976,  12646,   0, 10,  35181)
    22  = @("Alemania",                213468,  388,   194700,  9528,    0, 8,   9240)
    60  = @("Argelia",                 33055,   551,   22375,   9419,    0, 13,  1261)
    75  = @("Chequia",                 17387,   101,   11900,   5099,    0, 5,   388)
    98  = @("Albania",                 5889,    139,   3123,    2584,    0, 6,   182)
    99  = @("Paraguay",                5852,    0,     4645,    1148,    0, 0,   59)
    103 = @("Grecia",                  4973,    118,   1374,    3389,    0, 1,   210)
    120 = @("Cuba",                    2726,    25,    2396,    242,     0, 0,   88)
    137 = @("Tunez",                   1601,    17,    1233,    317,     0, 0,   51)
    150 = @("Siria",                   944,     52,    296,     600,     0, 2,   48)
    151 = @("Principado de Andorra",   939,     0,     825,     62,      0, 0,   52)
    152 = @("Republica del Chad",      938,     0,     814,     49,      0, 0,   75)
    153 = @("Malta",                   926,     36,    668,     249,     0, 0,   9)
    154 = @("Jamaica",                 920,     15,    745,     163,     0, 0,   12)
}

foreach ($r in $rows.Keys) {
    $entry = $rows[$r]
    $ws.Cells.Item($r, 1).Value = $entry[0]
    for ($col = 2; $col -le 8; $col++) {
        $ws.Cells.Item($r, $col).Value = $entry[$col - 1]
    }
}
